$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '30.331.54'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.933.55'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '0.7477'
$ws.Range("E5").Value = '  +3.18%  '
$ws.Range("D6").Value = '248.94'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '28.26'
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("D9").Value = '0.3209'
$ws.Range("E9").Value = '  -4.00%  '
$ws.Range("D10").Value = '0.07114'
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("D11").Value = '0.7890'
$ws.Range("E11").Value = '  -2.84%  '
$ws.Range("D12").Value = '0.08001'
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D13").Value = '1.931.87'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '5.389'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '94.56'
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").Value = '14.64'
$ws.Range("E16").Value = '  -2.40%  '
$ws.Range("D17").Value = '30.331.79'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '253.21'
$ws.Range("E18").Value = '  +1.42%  '
$ws.Range("D19").Value = '0.000008042'
$ws.Range("E19").Value = '  -2.56%  '
$ws.Range("D20").Value = '5.795'
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("D21").Value = '2.183.68'
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '6.818'
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("D25").Value = '9.575'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").Value = '164.57'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").Value = '2.340'
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").Value = '19.13'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").Value = '0.1329'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = '1.357'
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("D31").Value = '1.532'
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").Value = '4.450'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '4.150'
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("D34").Value = '0.05134'
$ws.Range("E34").Value = '  -1.39%  '
$ws.Range("D35").Value = '1.284'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").Value = '0.7502'
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").Value = '2.766'
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").Value = '0.01968'
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").Value = '2.806'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").Value = '78.23'
$ws.Range("E40").Value = '  -2.98%  '
$ws.Range("D41").Value = '6.419'
$ws.Range("E41").Value = '  -1.25%  '
$ws.Range("D42").Value = '0.4512'
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("D43").Value = '1.994'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").Value = '0.8436'
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '102.53'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = '9.827'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").Value = '7.540'
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").Value = '37.48'
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("D50").Value = '980.44'
$ws.Range("E50").Value = '  +11.71%  '
$ws.Range("D51").Value = '0.1199'
$ws.Range("E51").Value = '  +5.05%  '

$rng.ClearFormats()
